# Generate Report for Handback
# Renames the in-flight handback file from the 4cc044c5... GUID to the
# 13793503... GUID (refreshed timestamps) and appends a brand-new handback
# entry for a second file, db62ad9b-b343-4ade-b93f-b5ee8bc325e7.md, across
# the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "4cc044c5-93a8-4f3c-9be7-ec7391096bff"
$newGuid = "13793503-6295-4ec2-a89c-371c3f162bbe"
$addGuid = "db62ad9b-b343-4ade-b93f-b5ee8bc325e7"

$zhHash = "13d87bd803bd70fc72815d62d93ffa80b30ceaff"
$deHash = $zhHash
$addHash = "f25a62217b5c7b18cd395b58c1fcba6f71f03f97"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# NOTE: `Range.Hyperlinks.Delete()` in this host clears *every* hyperlink on
# the worksheet (it is not scoped to the range) -- so do a single clearing
# pass up front, then (re)add every hyperlink for the sheet in display order.
$ws1.Range("A1").Hyperlinks.Delete()

# -- row 2: rename the GUID, refresh the hyperlink + the generate date --
$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aea3cafb8ebaaed9bec5a5fc1600165fcf3e743e/e2e/$newGuid.md", "", "", "e2e\$newGuid.md")
$ws1.Range("C2").Value = ".md"
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("G2").Value = "2016-08-26 15:15:35"

# -- row 3: brand-new entry --
$ws1.Range("A3").Value = "$addGuid.md"
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aea3cafb8ebaaed9bec5a5fc1600165fcf3e743e/e2e/$addGuid.md", "", "", "e2e\$addGuid.md")
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText
$ws1.Range("G3").Value = "2016-08-26 15:15:35"

$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# NOTE: see comment above -- clear the whole sheet's hyperlinks once, then
# re-add every one of them (existing + new) in the right order.
$ws2.Range("A1").Hyperlinks.Delete()

# -- row 2: rename the GUID, refresh hyperlinks, xliff name + dates --
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aea3cafb8ebaaed9bec5a5fc1600165fcf3e743e/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = $statusText
$ws2.Range("D2").Value = "e2e"
$ws2.Range("E2").Value = "ht"
$ws2.Range("F2").Value = "'False"
$ws2.Range("G2").Value = "$newGuid.$zhHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-26 15:15:30"
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/37cbb7c059f240dac7ab69b52edda73ac1681f79/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws2.Range("J2").Value = "$newGuid.$zhHash.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-26 15:15:47"
$ws2.Range("L2").Value = "'"
$ws2.Range("M2").Value = "'True"
$ws2.Range("N2").Value = "'"
$ws2.Range("O2").Value = "'False"
$ws2.Range("P2").Value = "'"

# -- row 3: brand-new entry --
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aea3cafb8ebaaed9bec5a5fc1600165fcf3e743e/e2e/$addGuid.md", "", "", "$addGuid.md")
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = $statusText
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "'True"
$ws2.Range("G3").Value = "$addGuid.$addHash.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-26 15:15:30"
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/37cbb7c059f240dac7ab69b52edda73ac1681f79/e2e/$addGuid.md", "", "", "$addGuid.md")
$ws2.Range("J3").Value = "$addGuid.$addHash.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-26 15:15:47"
$ws2.Range("L3").Value = "'"
$ws2.Range("M3").Value = "'True"
$ws2.Range("N3").Value = "'"
$ws2.Range("O3").Value = "'False"
$ws2.Range("P3").Value = "'"

$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# NOTE: see comment above -- clear the whole sheet's hyperlinks once, then
# re-add every one of them (existing + new) in the right order.
$ws3.Range("A1").Hyperlinks.Delete()

# -- row 2: rename the GUID, refresh hyperlinks, xliff name + dates --
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aea3cafb8ebaaed9bec5a5fc1600165fcf3e743e/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = $statusText
$ws3.Range("D2").Value = "e2e"
$ws3.Range("E2").Value = "ht"
$ws3.Range("F2").Value = "'False"
$ws3.Range("G2").Value = "$newGuid.$deHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-26 15:15:35"
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/be825b93e90fe9331208c814a1c33ef117caba1b/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws3.Range("J2").Value = "$newGuid.$deHash.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-26 15:15:54"
$ws3.Range("L2").Value = "'"
$ws3.Range("M2").Value = "'True"
$ws3.Range("N2").Value = "'"
$ws3.Range("O2").Value = "'False"
$ws3.Range("P2").Value = "'"

# -- row 3: brand-new entry --
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aea3cafb8ebaaed9bec5a5fc1600165fcf3e743e/e2e/$addGuid.md", "", "", "$addGuid.md")
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = $statusText
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "'True"
$ws3.Range("G3").Value = "$addGuid.$addHash.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-26 15:15:35"
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/be825b93e90fe9331208c814a1c33ef117caba1b/e2e/$addGuid.md", "", "", "$addGuid.md")
$ws3.Range("J3").Value = "$addGuid.$addHash.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-26 15:15:54"
$ws3.Range("L3").Value = "'"
$ws3.Range("M3").Value = "'True"
$ws3.Range("N3").Value = "'"
$ws3.Range("O3").Value = "'False"
$ws3.Range("P3").Value = "'"

$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:P3"))
